# update 5/5 form + department
# Append a new "blank" department row (row 17) to the sheet, mirroring the
# layout/formatting of the existing rows, then leave the selection where the
# user last clicked while filling out the form (D12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 17

# New row data: Ma don vi / Ten don vi / Thuoc, all zero (placeholder values)
$ws.Cells.Item($newRow, 1).Value = 0
$ws.Cells.Item($newRow, 2).Value = 0
$ws.Cells.Item($newRow, 3).Value = 0

# Match the formatting used by column C's existing data rows
# (right-aligned, wrapped text - same as cell C16).
$ws.Cells.Item($newRow, 3).HorizontalAlignment = -4152
$ws.Cells.Item($newRow, 3).WrapText = $true

# Leave the active selection on D12, as in the saved workbook.
$ws.Range("D12").Select() | Out-Null
